# Fruta / hortaliza, semanal
# Insert a new weekly record as row 120 on the "Plátano" sheet, pushing the
# existing rows 120-138 down to 121-139 (dimension grows from T138 to T139).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 120..138 down by inserting a new blank row at 120.
$ws.Rows.Item(120).Insert()

# Populate the newly inserted row 120 with the new weekly observation.
$ws.Cells.Item(120, 1).Value = 1
$ws.Cells.Item(120, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(120, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(120, 4).Value = 44476
$ws.Cells.Item(120, 5).Value = 15
$ws.Cells.Item(120, 6).Value = "Fruta"
$ws.Cells.Item(120, 7).Value = 100108
$ws.Cells.Item(120, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(120, 9).Value = 100108006
$ws.Cells.Item(120, 10).Value = "Plátano"
$ws.Cells.Item(120, 11).Value = "Sin especificar"
$ws.Cells.Item(120, 12).Value = "Pintón"
$ws.Cells.Item(120, 13).Value = 120
$ws.Cells.Item(120, 14).Value = 20000
$ws.Cells.Item(120, 15).Value = 21000
$ws.Cells.Item(120, 16).Value = 20500
$ws.Cells.Item(120, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(120, 18).Value = "Ecuador"
$ws.Cells.Item(120, 19).Value = 1025
$ws.Cells.Item(120, 20).Value = 20
